$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.488.07"
$ws.Range("E2").Value = "  +0.55%  "

$ws.Range("D3").Value = "1.725.18"
$ws.Range("E3").Value = "  +0.14%  "

$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  +0.15%  "

$ws.Range("D5").Value = "'244.55"
$ws.Range("E5").Value = "  +2.05%  "

$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  +0.08%  "

$ws.Range("D7").Value = "'0.4800"
$ws.Range("E7").Value = "  +1.49%  "

$ws.Range("D8").Value = "'0.2666"
$ws.Range("E8").Value = "  +1.30%  "

$ws.Range("D9").Value = "'0.06208"
$ws.Range("E9").Value = "  -0.40%  "

$ws.Range("D10").Value = "1.725.73"
$ws.Range("E10").Value = "  +0.19%  "

$ws.Range("D11").Value = "'0.07155"
$ws.Range("E11").Value = "  +1.17%  "

$ws.Range("D12").Value = "'15.61"
$ws.Range("E12").Value = "  +2.05%  "

$ws.Range("D13").Value = "'0.6144"
$ws.Range("E13").Value = "  +3.55%  "

$ws.Range("D14").Value = "'4.512"
$ws.Range("E14").Value = "  +2.32%  "

$ws.Range("D15").Value = "'77.13"
$ws.Range("E15").Value = "  +1.00%  "

$ws.Range("D16").Value = "'1.000"
$ws.Range("E16").Value = "  +0.05%  "

$ws.Range("D17").Value = "26.506.61"
$ws.Range("E17").Value = "  +0.65%  "

$ws.Range("D18").Value = "'1.001"
$ws.Range("E18").Value = "  +0.12%  "

$ws.Range("D19").Value = "'0.000006916"
$ws.Range("E19").Value = "  +1.55%  "

$ws.Range("D20").Value = "'11.64"
$ws.Range("E20").Value = "  +0.28%  "

$ws.Range("D21").Value = "1.948.05"
$ws.Range("E21").Value = "  +0.57%  "

$ws.Range("D22").Value = "'4.518"
$ws.Range("E22").Value = "  -1.09%  "

$ws.Range("D23").Value = "'8.932"
$ws.Range("E23").Value = "  +1.68%  "

$ws.Range("D24").Value = "'5.261"
$ws.Range("E24").Value = "  -1.36%  "

$ws.Range("D25").Value = "'136.04"
$ws.Range("E25").Value = "  +0.79%  "

$ws.Range("D26").Value = "'15.36"
$ws.Range("E26").Value = "  +0.71%  "

$ws.Range("D27").Value = "'1.790"
$ws.Range("E27").Value = "  +1.49%  "

$ws.Range("D28").Value = "'1.405"
$ws.Range("E28").Value = "  -0.21%  "

$ws.Range("D29").Value = "'106.86"
$ws.Range("E29").Value = "  -0.08%  "

$ws.Range("D30").Value = "'3.968"
$ws.Range("E30").Value = "  -1.34%  "

$ws.Range("D31").Value = "'0.08020"
$ws.Range("E31").Value = "  +3.53%  "

$ws.Range("D32").Value = "'3.697"
$ws.Range("E32").Value = "  -0.11%  "

$ws.Range("D33").Value = "'0.04559"
$ws.Range("E33").Value = "  +2.44%  "

$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").Value = "'2.614"
$ws.Range("E34").Value = "  +0.09%  "

$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "'0.6351"
$ws.Range("E35").Value = "  +2.32%  "

$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").Value = "'0.9879"
$ws.Range("E36").Value = "  +1.19%  "

$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").Value = "'0.9292"
$ws.Range("E37").Value = "  +0.27%  "

$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").Value = "'2.083"
$ws.Range("E38").Value = "  +9.09%  "

$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").Value = "'2.420"
$ws.Range("E39").Value = "  -0.12%  "

$ws.Range("B40").Value = "Quant"
$ws.Range("C40").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D40").Value = "'104.70"
$ws.Range("E40").Value = "  -9.23%  "

$ws.Range("B41").Value = "PaxDollar"
$ws.Range("C41").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D41").Value = "'1.003"
$ws.Range("E41").Value = "  +0.25%  "

$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").Value = "'0.01500"
$ws.Range("E42").Value = "  +1.99%  "

$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "'5.612"
$ws.Range("E43").Value = "  +5.02%  "

$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").Value = "'0.3900"
$ws.Range("E44").Value = "  +2.03%  "

$ws.Range("B45").Value = "Aptos"
$ws.Range("C45").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D45").Value = "'6.900"
$ws.Range("E45").Value = "  +10.23%  "

$ws.Range("B46").Value = "Algorand"
$ws.Range("C46").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D46").Value = "'0.1186"
$ws.Range("E46").Value = "  +2.23%  "

$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").Value = "'0.05330"
$ws.Range("E47").Value = "  +0.74%  "

$ws.Range("B48").Value = "Elrond"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D48").Value = "'30.83"
$ws.Range("E48").Value = "  +0.75%  "

$ws.Range("D49").Value = "'7.854"
$ws.Range("E49").Value = "  +2.47%  "

$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").Value = "'1.266"
$ws.Range("E50").Value = "  +3.73%  "

$ws.Range("B51").Value = "Decentraland"
$ws.Range("C51").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D51").Value = "'0.3425"
$ws.Range("E51").Value = "  +0.95%  "
